# Update "想去人数" (F column) values for rows that changed in the
# refreshed data pull, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 7665
    9  = 5838
    12 = 25
    13 = 1789
    14 = 1298
    15 = 280
    16 = 197
    17 = 24
    18 = 5522
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
